$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) "Postgres Database" heading: remove the <w:lastRenderedPageBreak/> that
#    used to precede the heading text.
# ---------------------------------------------------------------------------
$postgresHeading = $d.Paragraphs.Item(54)
$r1 = $postgresHeading.Range.Duplicate
$xml1 = "<w:p $wns><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr><w:r><w:t>Postgres Database</w:t></w:r></w:p>"
$r1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) "Graphql/Relay" body paragraph:
#    - merge the two runs that make up the opening sentence into one run
#    - split "with this, Relay is a client..." into "with this, " and a new
#      run carrying <w:lastRenderedPageBreak/> before "Relay is a client..."
# ---------------------------------------------------------------------------
$graphqlBody = $d.Paragraphs.Item(52)
$r2 = $graphqlBody.Range.Duplicate
$xml2 = "<w:p $wns>" +
    "<w:r><w:t xml:space=`"preserve`">Graphql is a query language used for querying structured data from a backend. It is used instead of REST for making API calls. </w:t></w:r>" +
    "<w:r><w:t>Its</w:t></w:r>" +
    "<w:r><w:t xml:space=`"preserve`"> benefits </w:t></w:r>" +
    "<w:r><w:t>include</w:t></w:r>" +
    "<w:r><w:t xml:space=`"preserve`"> single REST Endpoint for all queries, requests reflect the structure of the returned data, a strict type system. Along </w:t></w:r>" +
    "<w:r><w:t xml:space=`"preserve`">with this, </w:t></w:r>" +
    "<w:r><w:lastRenderedPageBreak/><w:t>Relay is a client for graphql API’s, created by Facebook; that integrates directly with react and conforms to the flow architecture</w:t></w:r>" +
    "<w:r><w:t>.</w:t></w:r>" +
    "</w:p>"
$r2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) Insert a new "React-bootstrap" section (Heading2 + body paragraph +
#    trailing blank paragraph) right before the existing "Nestjs" heading,
#    without disturbing the blank paragraph that already precedes it.
# ---------------------------------------------------------------------------
$nestHeading = $d.Paragraphs.Item(46)
$insertPos = $nestHeading.Range.Start
$r3 = $d.Range($insertPos, $insertPos)
$xml3 = "<w:p $wns><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr><w:r><w:t>React-bootstrap</w:t></w:r></w:p>" +
    "<w:p $wns>" +
        "<w:r><w:t xml:space=`"preserve`">React-bootstrap is a </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:t>UI</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`"> component library for react. It will allow quick UI iteration, whilst providing a consistent style. I chose this component library as it has a basic layout system, easy to understand component documentation, and simplistic styling.</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wns/>"
$r3.InsertXML($xml3)
